$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header labels A10/A11/A12 pick up the bold "mtitleStyle" look (s=4) ---
# Copy the formatting already used by A9 (which carries that exact style)
# onto A10, A11 and A12 without touching their text values.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("A12").PasteSpecial(-4122)

# --- Score totals doubled ---
$ws.Range("D10").Value = 56
$ws.Range("E10").Value = 56

# --- Negative marking changed from -3 to -1 (keep it text, like before) ---
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "-1"
$ws.Range("C12").Copy()
$ws.Range("C11").PasteSpecial(-4122)

# --- New "Student Ans" / "Correct Ans" header columns G & H ---
$ws.Range("G15").Value = "Student Ans"
$ws.Range("A15").Copy()
$ws.Range("G15").PasteSpecial(-4122)

$ws.Range("H15").Value = "Correct Ans"
$ws.Range("A15").Copy()
$ws.Range("H15").PasteSpecial(-4122)

# --- Per-question "Correct Ans" values, columns D/E (rows 19-40) ---
$correctDE = @{
  19 = "Option A"; 20 = "Option D"; 21 = "Option B"; 22 = "Option C"; 23 = "Option B";
  24 = "Option C"; 25 = "Option D"; 26 = "Option D"; 27 = "Option A"; 28 = "Option A";
  29 = "Option C"; 30 = "Option A"; 31 = "Option D"; 32 = "Option D"; 33 = "Option B";
  34 = "Option D"; 35 = "Option C"; 36 = "Option D"; 37 = "Option B"; 38 = "Option D";
  39 = "Option A"; 40 = "Option A"
}

foreach ($r in 19..40) {
  $ws.Range("D16").Copy()
  $ws.Range("D$r").PasteSpecial(-4122)

  $ws.Range("E$r").Value = $correctDE[$r]
  $ws.Range("E16").Copy()
  $ws.Range("E$r").PasteSpecial(-4122)
}

# --- New "Student Ans" / "Correct Ans" values, columns G/H (rows 16-21) ---
$correctGH = @{
  16 = "Option A"; 17 = "Option D"; 18 = "Option D"; 19 = "Option A"; 20 = "Option C"; 21 = "Option D"
}

foreach ($r in 16..21) {
  $ws.Range("D16").Copy()
  $ws.Range("G$r").PasteSpecial(-4122)

  $ws.Range("H$r").Value = $correctGH[$r]
  $ws.Range("E16").Copy()
  $ws.Range("H$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
